$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update task descriptions to note tests are included
$ws.Range("B4").Value = "Create an issue class w/ tests"
$ws.Range("B5").Value = "Create a series class that holds issue items w/ tests"
$ws.Range("B6").Value = "Create a collection class that holds series items w/ tests"

# Mark these tasks as completed in Week 1 (Amount Remaining column)
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 1

# Update the active selection as left by the editing session
$ws.Range("E8").Select()
